{"js": "// Apply the documented text edits by paragraph index.\n// The body paragraphs are addressed positionally because several of\n// them share identical placeholder text (\"1\") that would otherwise be\n// ambiguous to a plain text search/replace.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst replacements = {\n  1: \"Introduction to cryptography\",\n  2: \"The large distribution is part of the most important sectors at the economic level for the French territory and a fortiori for the whole world. In 2020, there are more than 44000 food sales outlets listed and at least 10900 communes in France have at least one general grocery store.\",\n  3: \"The main brands present on the territory are Leclerc, Auchan, U, Intermarch\u00e9, Carrefour and Casino. They alone hold 85% of the market share. Hard discounters are also part of this landscape, with an increase in customer numbers in these stores of more than 60% over a period of 10 years. The sector's turnover exceeds \u20ac210 billion, with more than 700,000 employees, 89% of whom are directly employed on permanent contracts. The wage bill was \u20ac15.7 billion in 2020.\",\n  4: \"1) PESTEL analysis of the retail sector\",\n  5: \"'Politique' is translated as 'Policy'.\",\n  6: \"The political domain engages the stability of States and the ease of exchanging with the international for all that concerns certain types of products. Indeed, large distributions work regularly with foreign countries to be able to have in their rays a large number of references and thus satisfy consumers.\",\n  7: \"Note that each country has its own rules and the distributions must take account of various laws and regulations before bringing a product especially.\",\n  8: \"'Economic'\",\n  9: \"The large-scale retail trade has a great influence on the economy of a country because food products are essential for all consumers.\",\n  10: \"It is a sector where competition is certainly very high, but which also attracts a large number of customers every day, regardless of the brands. The health crisis did not affect this type of business, which remained open since it was obviously considered an essential business.\",\n  11: \"If customers tend to desert for some of them the too large hypermarkets in favor of smaller sales spaces, it remains that the brands of the sector do not know particular difficulties.\",\n  12: \"sociological\",\n  15: \"Most of them are loyal to customers through a loyalty card that allows them to accumulate points and earn discounts on all kinds of products.\",\n  16: \"Technological\",\n  17: \"As in all other sectors of activity, large-scale distribution has been forced to adapt to technological innovations. Automated checkouts or click and collect, everything is done to facilitate shopping and customer travel.\",\n  19: \"Hypermarkets are increasingly competing with specialized stores by highlighting sales areas dedicated to technology, with specially present sales advisors to help customers in their choices.\",\n  20: \"The large distribution wants to be in 2020 closer to the needs and expectations of consumers, with an increased presence on the net and especially on social networks. These last allow to fight more effectively against competitors by keeping a permanent interaction with Internet users.\",\n  22: \"The large distribution, like a majority of other companies, has made these last years a big effort to adapt itself to the ecological demand of partners as well as customers. Indeed, the organic products appeared in the shelves there are some years and references do not cease increasing still today. The consumers are very in demand on this type of product which little by little, extended to all the fields of activity, food certainly, but also household products and textile.\",\n  23: \"In addition, the bags become reusable at will by the customer and many brands deliver cardboard boxes for click and collect.\",\n  25: \"At the legal level, large distributions are subject to laws that protect consumers. Opening hours and days are regulated in order not to operate too unfair competition against other small businesses.\",\n  28: \"REFERENCES\",\n  29: \"According to Nielsen, sales in the large distribution sector increased by +2.5% in 2020 compared to 2019. This increase is mainly due to the health crisis and the closure of restaurants which has led consumers to buy more food products from supermarkets. In addition, online sales have also increased significantly (+20%) during this period.\",\n};\n\nfor (const key of Object.keys(replacements)) {\n  const idx = parseInt(key, 10);\n  paragraphs.items[idx].insertText(replacements[key], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the documented text edits by paragraph index (1-based, COM style).\n# Paragraphs are addressed positionally because several of them share\n# identical placeholder text (\"1\") that would otherwise be ambiguous to\n# a plain text search/replace.\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    2  = \"Introduction to cryptography\"\n    3  = \"The large distribution is part of the most important sectors at the economic level for the French territory and a fortiori for the whole world. In 2020, there are more than 44000 food sales outlets listed and at least 10900 communes in France have at least one general grocery store.\"\n    4  = \"The main brands present on the territory are Leclerc, Auchan, U, Intermarch\u00e9, Carrefour and Casino. They alone hold 85% of the market share. Hard discounters are also part of this landscape, with an increase in customer numbers in these stores of more than 60% over a period of 10 years. The sector's turnover exceeds \u20ac210 billion, with more than 700,000 employees, 89% of whom are directly employed on permanent contracts. The wage bill was \u20ac15.7 billion in 2020.\"\n    5  = \"1) PESTEL analysis of the retail sector\"\n    6  = \"'Politique' is translated as 'Policy'.\"\n    7  = \"The political domain engages the stability of States and the ease of exchanging with the international for all that concerns certain types of products. Indeed, large distributions work regularly with foreign countries to be able to have in their rays a large number of references and thus satisfy consumers.\"\n    8  = \"Note that each country has its own rules and the distributions must take account of various laws and regulations before bringing a product especially.\"\n    9  = \"'Economic'\"\n    10 = \"The large-scale retail trade has a great influence on the economy of a country because food products are essential for all consumers.\"\n    11 = \"It is a sector where competition is certainly very high, but which also attracts a large number of customers every day, regardless of the brands. The health crisis did not affect this type of business, which remained open since it was obviously considered an essential business.\"\n    12 = \"If customers tend to desert for some of them the too large hypermarkets in favor of smaller sales spaces, it remains that the brands of the sector do not know particular difficulties.\"\n    13 = \"sociological\"\n    16 = \"Most of them are loyal to customers through a loyalty card that allows them to accumulate points and earn discounts on all kinds of products.\"\n    17 = \"Technological\"\n    18 = \"As in all other sectors of activity, large-scale distribution has been forced to adapt to technological innovations. Automated checkouts or click and collect, everything is done to facilitate shopping and customer travel.\"\n    20 = \"Hypermarkets are increasingly competing with specialized stores by highlighting sales areas dedicated to technology, with specially present sales advisors to help customers in their choices.\"\n    21 = \"The large distribution wants to be in 2020 closer to the needs and expectations of consumers, with an increased presence on the net and especially on social networks. These last allow to fight more effectively against competitors by keeping a permanent interaction with Internet users.\"\n    23 = \"The large distribution, like a majority of other companies, has made these last years a big effort to adapt itself to the ecological demand of partners as well as customers. Indeed, the organic products appeared in the shelves there are some years and references do not cease increasing still today. The consumers are very in demand on this type of product which little by little, extended to all the fields of activity, food certainly, but also household products and textile.\"\n    24 = \"In addition, the bags become reusable at will by the customer and many brands deliver cardboard boxes for click and collect.\"\n    26 = \"At the legal level, large distributions are subject to laws that protect consumers. Opening hours and days are regulated in order not to operate too unfair competition against other small businesses.\"\n    29 = \"REFERENCES\"\n    30 = \"According to Nielsen, sales in the large distribution sector increased by +2.5% in 2020 compared to 2019. This increase is mainly due to the health crisis and the closure of restaurants which has led consumers to buy more food products from supermarkets. In addition, online sales have also increased significantly (+20%) during this period.\"\n}\n\nforeach ($idx in $replacements.Keys) {\n    $d.Paragraphs.Item($idx).Range.Text = $replacements[$idx]\n}\n"}
